$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "Issue date" text (row 5) ---
$ws.Range("A5").Value = "Issue date: 22/06/2021 11:46:03"

# --- Update existing result row (row 16) ---
$ws.Range("A16").Value = 669
$ws.Range("C16").Value = "Stylistic Features: slf, frc, e50te, agf, sxf, caf, anf, spe, nw, pnf, vof, acf, inf, aof, pw, thf"
$ws.Range("D16").Value = "None"
$ws.Range("F16").Value = "5 folds X 20 iterations CV"

# J16's new text ("91.73") looks numeric, so Excel would otherwise store it
# as a number. Force it to be kept as text (matching the original file,
# where it is a shared string), then restore the cell's original look
# (style 6, same as F11) since forcing text format mutates the style.
$ws.Range("J16").NumberFormat = "@"
$ws.Range("J16").Value = "91.73"
$ws.Range("F11").Copy()
$ws.Range("J16").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Add the new result row (row 17) ---
$ws.Range("A17").Value = 501
$ws.Range("B17").Value = "Hebrew"
$ws.Range("C17").Value = "Stylistic Features: vof, huf, aof, pnf, anf, agf, frc, mef, acf, fdf"
$ws.Range("D17").Value = "None"
$ws.Range("E17").Value = "lowercase"
$ws.Range("F17").Value = "5 folds X 20 iterations CV"
$ws.Range("J17").Value = "90.38V"

# Copy the formatting of row 16's plain "data" cells onto the new row 17
# cells (A:F, and J17 too - unlike J16, J17 uses the plain style, not the
# highlighted one).
$ws.Range("A16:F16").Copy()
$ws.Range("A17:F17").PasteSpecial(-4122)
$ws.Range("A16").Copy()
$ws.Range("J17").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Resize the table to include the new row ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A15:L17"))
$lo.TableStyle = "TableStyleLight9"

# --- Column width changes ---
# (ColumnWidth is quantized internally to 1/6-character steps by this
# engine, so 97.835 / 25.835 are the closest inputs that land on the
# target stored widths of ~98.71 / ~26.71 characters.)
$ws.Columns.Item(3).ColumnWidth = 97.835
$ws.Columns.Item(6).ColumnWidth = 25.835
